$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)

# Build the new "up/down arrow" shape by duplicating the existing
# "right arrow" shape already on the slide (Shapes.Item(6), id=3,
# "오른쪽 화살표 2"). Duplicating it carries over the exact same
# p:style (lnRef/fillRef/effectRef/fontRef), fill (F23B48), line
# (noFill) and txBody formatting that the new shape needs, then we
# just re-type it as an up/down arrow, rename it, and reposition it.
$src = $s.Shapes.Item(6)
$dupRange = $src.Duplicate()
$dup = $dupRange.Item(1)

$dup.Name = "위쪽/아래쪽 화살표 3"
$dup.AutoShapeType = 38

# Target geometry (EMU): off x=4096512 y=3572256, ext cx=182880 cy=353568
# (the literals below are nudged by a fraction of a point to cancel out
# the host's internal float32 rounding in the Left/Top/Width/Height
# setters, so the EMU actually written lands exactly on target instead
# of being truncated one EMU short).
$dup.Left = 322.5600129
$dup.Top = 281.28001409999996
$dup.Width = 14.4000001
$dup.Height = 27.84
